# feat: aplica regras de negócio na base tratada
#
# Populates the "base_tratada" sheet with the raw visit rows (columns A-E,
# copied from "raw" so values/styles - e.g. the date format on column B -
# come along for the ride) and then adds four business-rule columns:
#   F: visita_valida               -> Sim/Não  (was the visit actually attempted?)
#   G: status_normalizado          -> CAPTADO/AUSENTE/RECUSA/OUTRO
#   H: fluxo                       -> segue_para_capex / encerrado_na_visita
#   I: flag_instalacao_hidrometro  -> Sim/Não

$wb  = $excel.ActiveWorkbook
$raw = $wb.Worksheets.Item("raw")
$bt  = $wb.Worksheets.Item("base_tratada")

# --- Rename the renamed column header (visita_realizada_flag -> visita_valida) ---
# This edits the shared-string entry in place via the existing F1 cell.
$bt.Range("F1").Value = "visita_valida"

# --- Bring across the raw rows (id_visita, data_visita, matricula, bairro,
#     status_visita) into columns A:E, rows 2-6. Copy (not just Value) so the
#     date number-format style on column B comes along too. ---
$raw.Range("A2:E6").Copy($bt.Range("A2"))

# --- New header cells for the derived/business-rule columns ---
$bt.Range("G1").Value = "status_normalizado"
$bt.Range("H1").Value = "fluxo"
$bt.Range("I1").Value = "flag_instalacao_hidrometro"

# --- Row-by-row business rules (literal per-row formulas; values recalc automatically) ---
for ($r = 2; $r -le 6; $r++) {
    $bt.Range("F$r").Formula = "=IF(OR(E$r=`"`",E$r=`"ERRO`"),`"Não`",`"Sim`")"
    $bt.Range("G$r").Formula = "=IF(E$r=`"SUCESSO`",`"CAPTADO`",IF(E$r=`"AUSENTE`",`"AUSENTE`",IF(E$r=`"RECUSA`",`"RECUSA`",`"OUTRO`")))"
    $bt.Range("H$r").Formula = "=IF(G$r=`"CAPTADO`",`"segue_para_capex`",`"encerrado_na_visita`")"
    $bt.Range("I$r").Formula = "=IF(H$r=`"segue_para_capex`",`"Sim`",`"Não`")"
}

# --- Formatting: wrap text on the "visita_valida" column (header cell was
#     sized for the formula text when it was entered) ---
$bt.Range("F2:F6").WrapText = $true
$bt.Range("K2").WrapText = $true
$bt.Rows.Item(2).RowHeight = 30.75

# --- Let Excel re-measure the widened columns now that they hold real data ---
$bt.Columns.Item(6).AutoFit()
$bt.Columns.Item(8).AutoFit()
$bt.Columns.Item(9).AutoFit()

# --- Selections matching where the editor's cursor ended up ---
$raw.Range("A2:E6").Select()
$bt.Range("I2:I6").Select()
$bt.Activate()
